# Include the image's markdown title (in addition to the link, which was
# already present) in PowerPoint's description (alt text) of the image.
#
# Source markdown: ![alt text](lalune.jpg "fig:")
# -> pandoc encodes the link target into the `descr` attribute, and now also
#    folds the title ("fig:") in front of it, separated by two spaces:
#    descr="fig:  lalune.jpg"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)

$shp.AlternativeText = "fig:  " + $shp.AlternativeText
